# "feat: add detail jawaban page"
# Replace the five placeholder question rows (previously long Arabic verses
# repeated across every option column) with short example/test values:
#   B = "Soal ke N", C = "opsi A.N", D = "opsi B.N", E = "opsi C.N",
#   F = "opsi D.N", G = "opsi E.N"
# and drop the wrap-text / fixed-row-height formatting those cells had.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# First element of each row is the worksheet row number (column A / bobot
# is untouched by this edit); the rest are the new B..G values.
$data = @(
    @(3, "Soal ke 1", "opsi A.1", "opsi B.1", "opsi C.1", "opsi D.1", "opsi E.1"),
    @(4, "Soal ke 2", "opsi A.2", "opsi B.2", "opsi C.2", "opsi D.2", "opsi E.2"),
    @(5, "Soal ke 3", "opsi A.3", "opsi B.3", "opsi C.3", "opsi D.3", "opsi E.3"),
    @(6, "Soal ke 4", "opsi A.4", "opsi B.4", "opsi C.4", "opsi D.4", "opsi E.4"),
    @(7, "Soal ke 5", "opsi A.5", "opsi B.5", "opsi C.5", "opsi D.5", "opsi E.5")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# The old wrapText alignment style is no longer needed on these cells.
$ws.Range("B3:G7").Style = "Normal"

# Rows no longer need the explicit (wrapped-text driven) heights either.
$ws.Rows("3:7").AutoFit()

# Move the current selection on the frozen (bottom-left) pane to C10.
[void]$ws.Range("C10").Select()
